$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be treated as text while we write the new
# values (several of them are plain decimals like "1.008" that Excel
# would otherwise auto-convert to numbers), then drop the formatting
# override again so the cells end up back at the default (unstyled) xf,
# matching the original inline-string cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.165.02"
$ws.Range("D3").Value = "1.780.34"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.62%  "
$ws.Range("D5").Value = "334.34"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("D7").Value = "0.3770"
$ws.Range("E7").Value = "  -0.69%  "
$ws.Range("D8").Value = "0.3404"
$ws.Range("E8").Value = "  -1.83%  "
$ws.Range("D9").Value = "48.10"
$ws.Range("E9").Value = "  -2.46%  "
$ws.Range("D10").Value = "1.186"
$ws.Range("E10").Value = "  -2.44%  "
$ws.Range("D11").Value = "0.07407"
$ws.Range("E11").Value = "  -3.23%  "
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "21.43"
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").Value = "6.373"
$ws.Range("E14").Value = "  -2.92%  "
$ws.Range("D15").Value = "1.779.80"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "7.005"
$ws.Range("E16").Value = "  -2.87%  "
$ws.Range("D17").Value = "0.00001081"
$ws.Range("E17").Value = "  -2.84%  "
$ws.Range("D18").Value = "0.06652"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").Value = "83.77"
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").Value = "6.513"
$ws.Range("D22").Value = "17.13"
$ws.Range("E22").Value = "  -2.41%  "
$ws.Range("D23").Value = "27.174.32"
$ws.Range("E23").Value = "  -0.83%  "
$ws.Range("D24").Value = "12.35"
$ws.Range("E24").Value = "  -5.73%  "
$ws.Range("D25").Value = "2.416"
$ws.Range("E25").Value = "  -2.22%  "
$ws.Range("D26").Value = "1.484"
$ws.Range("E26").Value = "  -1.78%  "
$ws.Range("D27").Value = "2.499"
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("D28").Value = "20.98"
$ws.Range("E28").Value = "  +4.20%  "
$ws.Range("D29").Value = "151.98"
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("D30").Value = "1.981.78"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("D31").Value = "132.29"
$ws.Range("E31").Value = "  -2.04%  "
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("D33").Value = "5.958"
$ws.Range("E33").Value = "  -4.60%  "
$ws.Range("D34").Value = "0.08578"
$ws.Range("E34").Value = "  -1.69%  "
$ws.Range("D35").Value = "12.95"
$ws.Range("E35").Value = "  -3.38%  "
$ws.Range("D36").Value = "1.647"
$ws.Range("E36").Value = "  -4.11%  "
$ws.Range("D37").Value = "5.366"
$ws.Range("E37").Value = "  -4.05%  "
$ws.Range("D38").Value = "0.6771"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("D39").Value = "0.06289"
$ws.Range("E39").Value = "  -2.57%  "
$ws.Range("D40").Value = "0.02323"
$ws.Range("E40").Value = "  -3.66%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "0.2167"
$ws.Range("E41").Value = "  -3.35%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "8.687"
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("D44").Value = "14.42"
$ws.Range("E44").Value = "  -2.13%  "
$ws.Range("D45").Value = "1.003"
$ws.Range("E45").Value = "  +0.70%  "
$ws.Range("D46").Value = "0.6309"
$ws.Range("E46").Value = "  -1.78%  "
$ws.Range("D47").Value = "3.833"
$ws.Range("E47").Value = "  -3.14%  "
$ws.Range("D48").Value = "2.104"
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("D49").Value = "128.14"
$ws.Range("E49").Value = "  -2.07%  "
$ws.Range("D50").Value = "0.07163"
$ws.Range("E50").Value = "  -2.46%  "
$ws.Range("D51").Value = "78.79"
$ws.Range("E51").Value = "  -1.39%  "

$ws.Range("D2:D51").ClearFormats()
